# Aggiornamento dati San Cesario SP al 23 agosto 2021
# Appends rows 344-357 (2021-08-10 .. 2021-08-23) to the daily series in
# columns A:D, mirroring the style already used by the existing date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44418, 0, 5, 76.14986293024673),
    @(44419, 0, 4, 60.91989034419738),
    @(44420, 3, 6, 91.37983551629607),
    @(44421, 0, 6, 91.37983551629607),
    @(44422, 1, 5, 76.14986293024673),
    @(44423, 1, 6, 91.37983551629607),
    @(44424, 2, 7, 106.6098081023454),
    @(44425, 0, 7, 106.6098081023454),
    @(44426, 1, 8, 121.8397806883948),
    @(44427, 0, 5, 76.14986293024673),
    @(44428, 0, 5, 76.14986293024673),
    @(44429, 0, 4, 60.91989034419738),
    @(44430, 3, 6, 91.37983551629607),
    @(44431, 0, 4, 60.91989034419738)
)

$startRow = 344
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Column A holds dates styled like the rest of the series (bordered, bold,
# centered, custom date number format) - copy that formatting down onto the
# newly added date cells so the appended rows match the existing ones.
$ws.Range("A343").Copy()
$ws.Range("A344:A357").PasteSpecial(-4122)
$excel.CutCopyMode = 0
